$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update existing L3 (2020 -> 2021) and add new M3 (2022),
# copying L3's number formatting onto the new M3 cell first.
$ws.Cells.Item(3, 12).Value = 2021

$ws.Cells.Item(3, 12).Copy() | Out-Null
$ws.Cells.Item(3, 13).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 13).Value = 2022

# Row 4: add new M4 (6.18), copying L4's formatting onto it.
$ws.Cells.Item(4, 12).Copy() | Out-Null
$ws.Cells.Item(4, 13).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 13).Value = 6.18

$excel.CutCopyMode = $false

# Move the active selection from M12 to M9, matching the saved view state.
$ws.Range("M9").Select() | Out-Null
